$d = $word.ActiveDocument
$d.Content.Find.Execute("by <<hearingType>>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<<hearingType>>", 2)
